$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C column) dates for rows 2 through 18 from 45183 to 45184
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}

# Update the hyperlink formulas in row 2 (columns S through Y) to include the link text argument
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/artfynd/A 30779-2023.xlsx, "A 30779-2023"")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/kartor/A 30779-2023.png", "A 30779-2023")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/knärot/A 30779-2023.png", "A 30779-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/klagomål/A 30779-2023.docx", "A 30779-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/klagomålsmail/A 30779-2023.docx", "A 30779-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/tillsyn/A 30779-2023.docx", "A 30779-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/tillsynsmail/A 30779-2023.docx", "A 30779-2023")'
